# Update crypto price/volume data per Wed Apr 12 19:54:10 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.905.93"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "1.904.46"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'319.35"
$ws.Range("E5").Value = "  -2.09%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.5035"
$ws.Range("E7").Value = "  -2.65%  "
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").Value = "'0.08245"
$ws.Range("E9").Value = "  -2.68%  "
$ws.Range("D10").Value = "'41.95"
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("E11").Value = "  -2.22%  "
$ws.Range("D12").Value = "'23.96"
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("D13").Value = "1.904.07"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").Value = "'6.346"
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("D15").Value = "'7.182"
$ws.Range("E15").Value = "  -2.27%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "'91.80"
$ws.Range("E17").Value = "  -3.35%  "
$ws.Range("D18").Value = "'0.00001091"
$ws.Range("E18").Value = "  -2.08%  "
$ws.Range("D19").Value = "'0.06489"
$ws.Range("E19").Value = "  -2.78%  "
$ws.Range("D20").Value = "'17.98"
$ws.Range("E20").Value = "  -1.87%  "
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").Value = "'5.932"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").Value = "29.940.40"
$ws.Range("E23").Value = "  -1.08%  "
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("D25").Value = "'2.189"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("D26").Value = "'21.97"
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("D27").Value = "2.124.37"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").Value = "'161.30"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "'2.258"
$ws.Range("E29").Value = "  -5.97%  "
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("D31").Value = "'1.121"
$ws.Range("E31").Value = "  +2.54%  "
$ws.Range("E32").Value = "  -2.58%  "
$ws.Range("E33").Value = "  -2.14%  "
$ws.Range("D34").Value = "'3.798"
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "'5.363"
$ws.Range("E35").Value = "  +2.34%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.02427"
$ws.Range("E36").Value = "  -2.85%  "
$ws.Range("D37").Value = "'0.06329"
$ws.Range("E37").Value = "  -3.86%  "
$ws.Range("D38").Value = "'0.2139"
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("E39").Value = "  -2.78%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'8.620"
$ws.Range("E40").Value = "  -2.06%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6396"
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("D42").Value = "'11.29"
$ws.Range("E42").Value = "  -4.94%  "
$ws.Range("D43").Value = "'1.202"
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("D44").Value = "'2.188"
$ws.Range("E44").Value = "  +6.33%  "
$ws.Range("D45").Value = "'13.22"
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("D46").Value = "'0.5985"
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("D47").Value = "'3.632"
$ws.Range("E47").Value = "  -2.20%  "
$ws.Range("D48").Value = "'122.20"
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("E49").Value = "  -3.12%  "
$ws.Range("D50").Value = "'78.29"
$ws.Range("E50").Value = "  -1.23%  "
$ws.Range("D51").Value = "'1.128"
$ws.Range("E51").Value = "  -3.16%  "
